# Auto-generated: apply market-price / profit recalculation updates
# across the 8 Leve worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2184.9092  # H19: 2442.4443 -> 2184.9092
$ws.Cells.Item(19, 9).Value = 2200  # I19: 2900 -> 2200
$ws.Cells.Item(19, 10).Value = 2179.25  # J19: 2311.7144 -> 2179.25
$ws.Cells.Item(19, 11).Value = 2200  # K19: 2900 -> 2200
$ws.Cells.Item(19, 12).Value = 2179.25  # L19: 2311.7144 -> 2179.25
$ws.Cells.Item(19, 13).Value = -2025  # M19: -2725 -> -2025
$ws.Cells.Item(19, 14).Value = -2529.25  # N19: -2661.7144 -> -2529.25
$ws.Cells.Item(70, 8).Value = 2773.077  # H70: 2782.1428 -> 2773.077
$ws.Cells.Item(70, 9).Value = 2500  # I70: 2400 -> 2500
$ws.Cells.Item(70, 10).Value = 2795.8333  # J70: 2886.3635 -> 2795.8333
$ws.Cells.Item(70, 11).Value = 7500  # K70: 7200 -> 7500
$ws.Cells.Item(70, 12).Value = 8387.499899999999  # L70: 8659.0905 -> 8387.499899999999
$ws.Cells.Item(70, 13).Value = -7230  # M70: -6930 -> -7230
$ws.Cells.Item(70, 14).Value = -8927.499899999999  # N70: -9199.0905 -> -8927.499899999999
$ws.Cells.Item(73, 8).Value = 2773.077  # H73: 2782.1428 -> 2773.077
$ws.Cells.Item(73, 9).Value = 2500  # I73: 2400 -> 2500
$ws.Cells.Item(73, 10).Value = 2795.8333  # J73: 2886.3635 -> 2795.8333
$ws.Cells.Item(73, 11).Value = 7500  # K73: 7200 -> 7500
$ws.Cells.Item(73, 12).Value = 8387.499899999999  # L73: 8659.0905 -> 8387.499899999999
$ws.Cells.Item(73, 13).Value = -6564  # M73: -6264 -> -6564
$ws.Cells.Item(73, 14).Value = -10259.4999  # N73: -10531.0905 -> -10259.4999
$ws.Cells.Item(80, 8).Value = 1060.3529  # H80: 1310.5834 -> 1060.3529
$ws.Cells.Item(80, 9).Value = 883.4  # I80: 1155.6666 -> 883.4
$ws.Cells.Item(80, 10).Value = 1313.1428  # J80: 1465.5 -> 1313.1428
$ws.Cells.Item(80, 11).Value = 2650.2  # K80: 3466.9998 -> 2650.2
$ws.Cells.Item(80, 12).Value = 3939.4284  # L80: 4396.5 -> 3939.4284
$ws.Cells.Item(80, 13).Value = -1652.2  # M80: -2468.9998 -> -1652.2
$ws.Cells.Item(80, 14).Value = -5935.428400000001  # N80: -6392.5 -> -5935.428400000001
$ws.Cells.Item(83, 8).Value = 1060.3529  # H83: 1310.5834 -> 1060.3529
$ws.Cells.Item(83, 9).Value = 883.4  # I83: 1155.6666 -> 883.4
$ws.Cells.Item(83, 10).Value = 1313.1428  # J83: 1465.5 -> 1313.1428
$ws.Cells.Item(83, 11).Value = 7950.599999999999  # K83: 10400.9994 -> 7950.599999999999
$ws.Cells.Item(83, 12).Value = 11818.2852  # L83: 13189.5 -> 11818.2852
$ws.Cells.Item(83, 13).Value = -2958.599999999999  # M83: -5408.999400000001 -> -2958.599999999999
$ws.Cells.Item(83, 14).Value = -21802.2852  # N83: -23173.5 -> -21802.2852
$ws.Cells.Item(106, 8).Value = 3189.75  # H106: 3123.524 -> 3189.75
$ws.Cells.Item(106, 9).Value = 2351.75  # I106: 2309.2307 -> 2351.75
$ws.Cells.Item(106, 11).Value = 2351.75  # K106: 2309.2307 -> 2351.75
$ws.Cells.Item(106, 13).Value = -1720.75  # M106: -1678.2307 -> -1720.75
$ws.Cells.Item(107, 8).Value = 1249.75  # H107: 1006 -> 1249.75
$ws.Cells.Item(107, 10).Value = 999.6667  # J107: 509 -> 999.6667
$ws.Cells.Item(107, 12).Value = 999.6667  # L107: 509 -> 999.6667
$ws.Cells.Item(107, 14).Value = -4839.6667  # N107: -4349 -> -4839.6667
$ws.Cells.Item(131, 8).Value = 7354.778  # H131: 6067.5454 -> 7354.778
$ws.Cells.Item(131, 9).Value = 5774.125  # I131: 4674.3 -> 5774.125
$ws.Cells.Item(131, 11).Value = 17322.375  # K131: 14022.9 -> 17322.375
$ws.Cells.Item(131, 13).Value = -12282.375  # M131: -8982.900000000001 -> -12282.375
$ws.Cells.Item(132, 8).Value = 1491.4186  # H132: 1512.3658 -> 1491.4186
$ws.Cells.Item(132, 9).Value = 1439.6389  # I132: 1461.8529 -> 1439.6389
$ws.Cells.Item(132, 11).Value = 4318.9167  # K132: 4385.5587 -> 4318.9167
$ws.Cells.Item(132, 13).Value = -1788.9167  # M132: -1855.5587 -> -1788.9167
$ws.Cells.Item(137, 8).Value = 1751.3636  # H137: 1696.6666 -> 1751.3636
$ws.Cells.Item(137, 9).Value = 1585.1111  # I137: 1536.1 -> 1585.1111
$ws.Cells.Item(137, 11).Value = 4755.3333  # K137: 4608.299999999999 -> 4755.3333
$ws.Cells.Item(137, 13).Value = -2205.3333  # M137: -2058.299999999999 -> -2205.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7289.3  # H32: 7593.763 -> 7289.3
$ws.Cells.Item(32, 9).Value = 6751.2256  # I32: 7113.069 -> 6751.2256
$ws.Cells.Item(32, 11).Value = 6751.2256  # K32: 7113.069 -> 6751.2256
$ws.Cells.Item(32, 13).Value = -6464.2256  # M32: -6826.069 -> -6464.2256
$ws.Cells.Item(45, 8).Value = 4350.3687  # H45: 3749.739 -> 4350.3687
$ws.Cells.Item(45, 9).Value = 3393.25  # I45: 2769.125 -> 3393.25
$ws.Cells.Item(45, 11).Value = 3393.25  # K45: 2769.125 -> 3393.25
$ws.Cells.Item(45, 13).Value = -3016.25  # M45: -2392.125 -> -3016.25
$ws.Cells.Item(74, 8).Value = 4197.641  # H74: 4403.108 -> 4197.641
$ws.Cells.Item(74, 9).Value = 2225.2812  # I74: 2347.2 -> 2225.2812
$ws.Cells.Item(74, 11).Value = 2225.2812  # K74: 2347.2 -> 2225.2812
$ws.Cells.Item(74, 13).Value = -1351.2812  # M74: -1473.2 -> -1351.2812
$ws.Cells.Item(77, 8).Value = 4197.641  # H77: 4403.108 -> 4197.641
$ws.Cells.Item(77, 9).Value = 2225.2812  # I77: 2347.2 -> 2225.2812
$ws.Cells.Item(77, 11).Value = 11126.406  # K77: 11736 -> 11126.406
$ws.Cells.Item(77, 13).Value = -6758.405999999999  # M77: -7368 -> -6758.405999999999
$ws.Cells.Item(88, 8).Value = 1993.4375  # H88: 1932 -> 1993.4375
$ws.Cells.Item(88, 10).Value = 1927.3334  # J88: 1829.5 -> 1927.3334
$ws.Cells.Item(88, 12).Value = 1927.3334  # L88: 1829.5 -> 1927.3334
$ws.Cells.Item(88, 14).Value = -2739.3334  # N88: -2641.5 -> -2739.3334
$ws.Cells.Item(91, 8).Value = 1993.4375  # H91: 1932 -> 1993.4375
$ws.Cells.Item(91, 10).Value = 1927.3334  # J91: 1829.5 -> 1927.3334
$ws.Cells.Item(91, 12).Value = 1927.3334  # L91: 1829.5 -> 1927.3334
$ws.Cells.Item(91, 14).Value = -4735.3334  # N91: -4637.5 -> -4735.3334
$ws.Cells.Item(102, 8).Value = 3330.8096  # H102: 3444.85 -> 3330.8096
$ws.Cells.Item(102, 9).Value = 2349.8235  # I102: 2431.0625 -> 2349.8235
$ws.Cells.Item(102, 11).Value = 2349.8235  # K102: 2431.0625 -> 2349.8235
$ws.Cells.Item(102, 13).Value = -727.8235  # M102: -809.0625 -> -727.8235
$ws.Cells.Item(110, 8).Value = 7749.8335  # H110: 7954.4546 -> 7749.8335
$ws.Cells.Item(110, 9).Value = 6555.3335  # I110: 6687.375 -> 6555.3335
$ws.Cells.Item(110, 11).Value = 6555.3335  # K110: 6687.375 -> 6555.3335
$ws.Cells.Item(110, 13).Value = -4510.3335  # M110: -4642.375 -> -4510.3335
$ws.Cells.Item(122, 8).Value = 1664.5625  # H122: 1788.3846 -> 1664.5625
$ws.Cells.Item(122, 9).Value = 1664.5625  # I122: 1788.3846 -> 1664.5625
$ws.Cells.Item(122, 11).Value = 4993.6875  # K122: 5365.1538 -> 4993.6875
$ws.Cells.Item(122, 13).Value = -2543.6875  # M122: -2915.1538 -> -2543.6875
$ws.Cells.Item(132, 8).Value = 1822.909  # H132: 1876.2903 -> 1822.909
$ws.Cells.Item(132, 9).Value = 1822.909  # I132: 1876.2903 -> 1822.909
$ws.Cells.Item(132, 11).Value = 5468.727000000001  # K132: 5628.8709 -> 5468.727000000001
$ws.Cells.Item(132, 13).Value = -2938.727000000001  # M132: -3098.8709 -> -2938.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2541.4482  # H20: 2483 -> 2541.4482
$ws.Cells.Item(20, 9).Value = 2204.2  # I20: 2115.6875 -> 2204.2
$ws.Cells.Item(20, 11).Value = 2204.2  # K20: 2115.6875 -> 2204.2
$ws.Cells.Item(20, 13).Value = -1957.2  # M20: -1868.6875 -> -1957.2
$ws.Cells.Item(81, 8).Value = 47500  # H81: 36696 -> 47500
$ws.Cells.Item(81, 10).Value = 47500  # J81: 36696 -> 47500
$ws.Cells.Item(81, 12).Value = 47500  # L81: 36696 -> 47500
$ws.Cells.Item(81, 14).Value = -49622  # N81: -38818 -> -49622
$ws.Cells.Item(84, 8).Value = 47500  # H84: 36696 -> 47500
$ws.Cells.Item(84, 10).Value = 47500  # J84: 36696 -> 47500
$ws.Cells.Item(84, 12).Value = 142500  # L84: 110088 -> 142500
$ws.Cells.Item(84, 14).Value = -153108  # N84: -120696 -> -153108
$ws.Cells.Item(86, 8).Value = 1346.909  # H86: 1233.9615 -> 1346.909
$ws.Cells.Item(86, 9).Value = 1146.1875  # I86: 1042.1578 -> 1146.1875
$ws.Cells.Item(86, 10).Value = 1882.1666  # J86: 1754.5714 -> 1882.1666
$ws.Cells.Item(86, 11).Value = 1146.1875  # K86: 1042.1578 -> 1146.1875
$ws.Cells.Item(86, 12).Value = 1882.1666  # L86: 1754.5714 -> 1882.1666
$ws.Cells.Item(86, 13).Value = -23.1875  # M86: 80.84220000000005 -> -23.1875
$ws.Cells.Item(86, 14).Value = -4128.1666  # N86: -4000.5714 -> -4128.1666
$ws.Cells.Item(89, 8).Value = 1346.909  # H89: 1233.9615 -> 1346.909
$ws.Cells.Item(89, 9).Value = 1146.1875  # I89: 1042.1578 -> 1146.1875
$ws.Cells.Item(89, 10).Value = 1882.1666  # J89: 1754.5714 -> 1882.1666
$ws.Cells.Item(89, 11).Value = 5730.9375  # K89: 5210.789 -> 5730.9375
$ws.Cells.Item(89, 12).Value = 9410.833000000001  # L89: 8772.857 -> 9410.833000000001
$ws.Cells.Item(89, 13).Value = -114.9375  # M89: 405.2110000000002 -> -114.9375
$ws.Cells.Item(89, 14).Value = -20642.833  # N89: -20004.857 -> -20642.833
$ws.Cells.Item(105, 8).Value = 2629.2222  # H105: 2695.625 -> 2629.2222
$ws.Cells.Item(105, 9).Value = 2629.2222  # I105: 2695.625 -> 2629.2222
$ws.Cells.Item(105, 11).Value = 2629.2222  # K105: 2695.625 -> 2629.2222
$ws.Cells.Item(105, 13).Value = -882.2222000000002  # M105: -948.625 -> -882.2222000000002
$ws.Cells.Item(132, 8).Value = 97499  # H132: 97374.5 -> 97499
$ws.Cells.Item(132, 10).Value = 97499  # J132: 97374.5 -> 97499
$ws.Cells.Item(132, 12).Value = 97499  # L132: 97374.5 -> 97499
$ws.Cells.Item(132, 14).Value = -107619  # N132: -107494.5 -> -107619
$ws.Cells.Item(134, 8).Value = 3852.238  # H134: 3813.7144 -> 3852.238
$ws.Cells.Item(134, 9).Value = 3734.6604  # I134: 3673.7407 -> 3734.6604
$ws.Cells.Item(134, 10).Value = 4475.4  # J134: 4653.5557 -> 4475.4
$ws.Cells.Item(134, 11).Value = 11203.9812  # K134: 11021.2221 -> 11203.9812
$ws.Cells.Item(134, 12).Value = 13426.2  # L134: 13960.6671 -> 13426.2
$ws.Cells.Item(134, 13).Value = -8668.9812  # M134: -8486.222099999999 -> -8668.9812
$ws.Cells.Item(134, 14).Value = -18496.2  # N134: -19030.6671 -> -18496.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4082  # H16: 4175.706 -> 4082
$ws.Cells.Item(16, 9).Value = 3419.8572  # I16: 3491.4614 -> 3419.8572
$ws.Cells.Item(16, 11).Value = 3419.8572  # K16: 3491.4614 -> 3419.8572
$ws.Cells.Item(16, 13).Value = -3132.8572  # M16: -3204.4614 -> -3132.8572
$ws.Cells.Item(31, 8).Value = 4267.0713  # H31: 4465.324 -> 4267.0713
$ws.Cells.Item(31, 10).Value = 4494.0967  # J31: 4819.885 -> 4494.0967
$ws.Cells.Item(31, 12).Value = 4494.0967  # L31: 4819.885 -> 4494.0967
$ws.Cells.Item(31, 14).Value = -5084.0967  # N31: -5409.885 -> -5084.0967
$ws.Cells.Item(34, 8).Value = 4267.0713  # H34: 4465.324 -> 4267.0713
$ws.Cells.Item(34, 10).Value = 4494.0967  # J34: 4819.885 -> 4494.0967
$ws.Cells.Item(34, 12).Value = 4494.0967  # L34: 4819.885 -> 4494.0967
$ws.Cells.Item(34, 14).Value = -4898.0967  # N34: -5223.885 -> -4898.0967
$ws.Cells.Item(105, 8).Value = 3358  # H105: 3695 -> 3358
$ws.Cells.Item(105, 9).Value = 3358  # I105: 3695 -> 3358
$ws.Cells.Item(105, 11).Value = 3358  # K105: 3695 -> 3358
$ws.Cells.Item(105, 13).Value = -1611  # M105: -1948 -> -1611
$ws.Cells.Item(113, 8).Value = 4082  # H113: 4175.706 -> 4082
$ws.Cells.Item(113, 9).Value = 3419.8572  # I113: 3491.4614 -> 3419.8572
$ws.Cells.Item(113, 11).Value = 3419.8572  # K113: 3491.4614 -> 3419.8572
$ws.Cells.Item(113, 13).Value = -1249.8572  # M113: -1321.4614 -> -1249.8572
$ws.Cells.Item(116, 8).Value = 70000  # H116: 56500 -> 70000
$ws.Cells.Item(116, 10).Value = 70000  # J116: 56500 -> 70000
$ws.Cells.Item(116, 12).Value = 70000  # L116: 56500 -> 70000
$ws.Cells.Item(116, 14).Value = -79178  # N116: -65678 -> -79178

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 184.46666  # H14: 191.21428 -> 184.46666
$ws.Cells.Item(14, 9).Value = 184.46666  # I14: 191.21428 -> 184.46666
$ws.Cells.Item(14, 11).Value = 553.3999799999999  # K14: 573.64284 -> 553.3999799999999
$ws.Cells.Item(14, 13).Value = -380.3999799999999  # M14: -400.64284 -> -380.3999799999999
$ws.Cells.Item(68, 8).Value = 45458476  # H68: 41670480 -> 45458476
$ws.Cells.Item(68, 9).Value = 71432536  # I68: 62503780 -> 71432536
$ws.Cells.Item(68, 11).Value = 214297608  # K68: 187511340 -> 214297608
$ws.Cells.Item(68, 13).Value = -214296797  # M68: -187510529 -> -214296797
$ws.Cells.Item(71, 8).Value = 45458476  # H71: 41670480 -> 45458476
$ws.Cells.Item(71, 9).Value = 71432536  # I71: 62503780 -> 71432536
$ws.Cells.Item(71, 11).Value = 642892824  # K71: 562534020 -> 642892824
$ws.Cells.Item(71, 13).Value = -642888768  # M71: -562529964 -> -642888768
$ws.Cells.Item(113, 8).Value = 448.83334  # H113: 465.66666 -> 448.83334
$ws.Cells.Item(113, 9).Value = 423.25  # I113: 448.5 -> 423.25
$ws.Cells.Item(113, 11).Value = 1269.75  # K113: 1345.5 -> 1269.75
$ws.Cells.Item(113, 13).Value = 900.25  # M113: 824.5 -> 900.25
$ws.Cells.Item(137, 8).Value = 5486.533  # H137: 5043.4707 -> 5486.533
$ws.Cells.Item(137, 9).Value = 3602.625  # I137: 3452 -> 3602.625
$ws.Cells.Item(137, 10).Value = 7639.5713  # J137: 7317 -> 7639.5713
$ws.Cells.Item(137, 11).Value = 10807.875  # K137: 10356 -> 10807.875
$ws.Cells.Item(137, 12).Value = 22918.7139  # L137: 21951 -> 22918.7139
$ws.Cells.Item(137, 13).Value = -5707.875  # M137: -5256 -> -5707.875
$ws.Cells.Item(137, 14).Value = -33118.7139  # N137: -32151 -> -33118.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 18000  # H70: 12333 -> 18000
$ws.Cells.Item(70, 9).Value = 0  # I70: 999 -> 0
$ws.Cells.Item(70, 11).Value = 0  # K70: 999 -> 0
$ws.Cells.Item(70, 13).ClearContents()  # M70: was -729
$ws.Cells.Item(73, 8).Value = 18000  # H73: 12333 -> 18000
$ws.Cells.Item(73, 9).Value = 0  # I73: 999 -> 0
$ws.Cells.Item(73, 11).Value = 0  # K73: 999 -> 0
$ws.Cells.Item(73, 13).ClearContents()  # M73: was -63
$ws.Cells.Item(80, 8).Value = 2608.8333  # H80: 2771.2273 -> 2608.8333
$ws.Cells.Item(80, 9).Value = 2418.5  # I80: 2610 -> 2418.5
$ws.Cells.Item(80, 10).Value = 2744.7856  # J80: 2882.8462 -> 2744.7856
$ws.Cells.Item(80, 11).Value = 2418.5  # K80: 2610 -> 2418.5
$ws.Cells.Item(80, 12).Value = 2744.7856  # L80: 2882.8462 -> 2744.7856
$ws.Cells.Item(80, 13).Value = -1420.5  # M80: -1612 -> -1420.5
$ws.Cells.Item(80, 14).Value = -4740.7856  # N80: -4878.8462 -> -4740.7856
$ws.Cells.Item(83, 8).Value = 2608.8333  # H83: 2771.2273 -> 2608.8333
$ws.Cells.Item(83, 9).Value = 2418.5  # I83: 2610 -> 2418.5
$ws.Cells.Item(83, 10).Value = 2744.7856  # J83: 2882.8462 -> 2744.7856
$ws.Cells.Item(83, 11).Value = 12092.5  # K83: 13050 -> 12092.5
$ws.Cells.Item(83, 12).Value = 13723.928  # L83: 14414.231 -> 13723.928
$ws.Cells.Item(83, 13).Value = -7100.5  # M83: -8058 -> -7100.5
$ws.Cells.Item(83, 14).Value = -23707.928  # N83: -24398.231 -> -23707.928
$ws.Cells.Item(97, 8).Value = 654.2973  # H97: 686.9722 -> 654.2973
$ws.Cells.Item(97, 9).Value = 606.4074000000001  # I97: 651.04 -> 606.4074000000001
$ws.Cells.Item(97, 10).Value = 783.6  # J97: 768.63635 -> 783.6
$ws.Cells.Item(97, 11).Value = 606.4074000000001  # K97: 651.04 -> 606.4074000000001
$ws.Cells.Item(97, 12).Value = 783.6  # L97: 768.63635 -> 783.6
$ws.Cells.Item(97, 13).Value = -110.4074000000001  # M97: -155.04 -> -110.4074000000001
$ws.Cells.Item(97, 14).Value = -1775.6  # N97: -1760.63635 -> -1775.6
$ws.Cells.Item(107, 8).Value = 1211.4117  # H107: 1249.625 -> 1211.4117
$ws.Cells.Item(107, 10).Value = 1188.2222  # J107: 1261.75 -> 1188.2222
$ws.Cells.Item(107, 12).Value = 1188.2222  # L107: 1261.75 -> 1188.2222
$ws.Cells.Item(107, 14).Value = -5028.2222  # N107: -5101.75 -> -5028.2222
$ws.Cells.Item(113, 8).Value = 2474.5715  # H113: 2638 -> 2474.5715
$ws.Cells.Item(113, 9).Value = 1497.25  # I113: 1498.3334 -> 1497.25
$ws.Cells.Item(113, 11).Value = 1497.25  # K113: 1498.3334 -> 1497.25
$ws.Cells.Item(113, 13).Value = 672.75  # M113: 671.6666 -> 672.75
$ws.Cells.Item(122, 8).Value = 2157.4211  # H122: 2158.4119 -> 2157.4211
$ws.Cells.Item(122, 9).Value = 2035.4286  # I122: 2053.6155 -> 2035.4286
$ws.Cells.Item(122, 11).Value = 6106.2858  # K122: 6160.8465 -> 6106.2858
$ws.Cells.Item(122, 13).Value = -3656.2858  # M122: -3710.8465 -> -3656.2858
$ws.Cells.Item(132, 8).Value = 4243.037  # H132: 4321.654 -> 4243.037
$ws.Cells.Item(132, 9).Value = 4312.048  # I132: 4312.095 -> 4312.048
$ws.Cells.Item(132, 10).Value = 4001.5  # J132: 4361.8 -> 4001.5
$ws.Cells.Item(132, 11).Value = 12936.144  # K132: 12936.285 -> 12936.144
$ws.Cells.Item(132, 12).Value = 12004.5  # L132: 13085.4 -> 12004.5
$ws.Cells.Item(132, 13).Value = -10406.144  # M132: -10406.285 -> -10406.144
$ws.Cells.Item(132, 14).Value = -17064.5  # N132: -18145.4 -> -17064.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3390.8333  # H40: 3469.9583 -> 3390.8333
$ws.Cells.Item(40, 9).Value = 3143.389  # I40: 3248.889 -> 3143.389
$ws.Cells.Item(40, 11).Value = 3143.389  # K40: 3248.889 -> 3143.389
$ws.Cells.Item(40, 13).Value = -3007.389  # M40: -3112.889 -> -3007.389
$ws.Cells.Item(45, 8).Value = 5000  # H45: 6000 -> 5000
$ws.Cells.Item(45, 9).Value = 5000  # I45: 6000 -> 5000
$ws.Cells.Item(45, 11).Value = 5000  # K45: 6000 -> 5000
$ws.Cells.Item(45, 13).Value = -4593  # M45: -5593 -> -4593
$ws.Cells.Item(60, 8).Value = 0  # H60: 48000 -> 0
$ws.Cells.Item(60, 10).Value = 0  # J60: 48000 -> 0
$ws.Cells.Item(60, 12).Value = 0  # L60: 48000 -> 0
$ws.Cells.Item(60, 14).ClearContents()  # N60: was -49018
$ws.Cells.Item(68, 8).Value = 5405.2  # H68: 5003.636 -> 5405.2
$ws.Cells.Item(68, 9).Value = 2024.8334  # I68: 1876.7142 -> 2024.8334
$ws.Cells.Item(68, 11).Value = 2024.8334  # K68: 1876.7142 -> 2024.8334
$ws.Cells.Item(68, 13).Value = -1275.8334  # M68: -1127.7142 -> -1275.8334
$ws.Cells.Item(71, 8).Value = 5405.2  # H71: 5003.636 -> 5405.2
$ws.Cells.Item(71, 9).Value = 2024.8334  # I71: 1876.7142 -> 2024.8334
$ws.Cells.Item(71, 11).Value = 10124.167  # K71: 9383.571 -> 10124.167
$ws.Cells.Item(71, 13).Value = -6380.166999999999  # M71: -5639.571 -> -6380.166999999999
$ws.Cells.Item(82, 8).Value = 991.0909  # H82: 1118.2 -> 991.0909
$ws.Cells.Item(82, 9).Value = 987.4286  # I82: 1154.4 -> 987.4286
$ws.Cells.Item(82, 10).Value = 997.5  # J82: 1082 -> 997.5
$ws.Cells.Item(82, 11).Value = 987.4286  # K82: 1154.4 -> 987.4286
$ws.Cells.Item(82, 12).Value = 997.5  # L82: 1082 -> 997.5
$ws.Cells.Item(82, 13).Value = -626.4286  # M82: -793.4000000000001 -> -626.4286
$ws.Cells.Item(82, 14).Value = -1719.5  # N82: -1804 -> -1719.5
$ws.Cells.Item(85, 8).Value = 991.0909  # H85: 1118.2 -> 991.0909
$ws.Cells.Item(85, 9).Value = 987.4286  # I85: 1154.4 -> 987.4286
$ws.Cells.Item(85, 10).Value = 997.5  # J85: 1082 -> 997.5
$ws.Cells.Item(85, 11).Value = 987.4286  # K85: 1154.4 -> 987.4286
$ws.Cells.Item(85, 12).Value = 997.5  # L85: 1082 -> 997.5
$ws.Cells.Item(85, 13).Value = 260.5714  # M85: 93.59999999999991 -> 260.5714
$ws.Cells.Item(85, 14).Value = -3493.5  # N85: -3578 -> -3493.5
$ws.Cells.Item(122, 8).Value = 2571.7273  # H122: 2580.818 -> 2571.7273
$ws.Cells.Item(122, 9).Value = 2521.111  # I122: 2558.9 -> 2521.111
$ws.Cells.Item(122, 10).Value = 2799.5  # J122: 2800 -> 2799.5
$ws.Cells.Item(122, 11).Value = 7563.333  # K122: 7676.700000000001 -> 7563.333
$ws.Cells.Item(122, 12).Value = 8398.5  # L122: 8400 -> 8398.5
$ws.Cells.Item(122, 13).Value = -5113.333  # M122: -5226.700000000001 -> -5113.333
$ws.Cells.Item(122, 14).Value = -13298.5  # N122: -13300 -> -13298.5
$ws.Cells.Item(136, 8).Value = 44449710  # H136: 46517096 -> 44449710
$ws.Cells.Item(136, 9).Value = 29417190  # I136: 30308590 -> 29417190
$ws.Cells.Item(136, 10).Value = 90913880  # J136: 100005170 -> 90913880
$ws.Cells.Item(136, 11).Value = 88251570  # K136: 90925770 -> 88251570
$ws.Cells.Item(136, 12).Value = 272741640  # L136: 300015510 -> 272741640
$ws.Cells.Item(136, 13).Value = -88249020  # M136: -90923220 -> -88249020
$ws.Cells.Item(136, 14).Value = -272746740  # N136: -300020610 -> -272746740

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3641.6667  # H62: 3553.8462 -> 3641.6667
$ws.Cells.Item(62, 10).Value = 3485.7144  # J62: 3362.5 -> 3485.7144
$ws.Cells.Item(62, 12).Value = 3485.7144  # L62: 3362.5 -> 3485.7144
$ws.Cells.Item(62, 14).Value = -4733.7144  # N62: -4610.5 -> -4733.7144
$ws.Cells.Item(65, 8).Value = 3641.6667  # H65: 3553.8462 -> 3641.6667
$ws.Cells.Item(65, 10).Value = 3485.7144  # J65: 3362.5 -> 3485.7144
$ws.Cells.Item(65, 12).Value = 17428.572  # L65: 16812.5 -> 17428.572
$ws.Cells.Item(65, 14).Value = -23668.572  # N65: -23052.5 -> -23668.572
$ws.Cells.Item(80, 8).Value = 0  # H80: 30000 -> 0
$ws.Cells.Item(80, 10).Value = 0  # J80: 30000 -> 0
$ws.Cells.Item(80, 12).Value = 0  # L80: 30000 -> 0
$ws.Cells.Item(80, 14).ClearContents()  # N80: was -31996
$ws.Cells.Item(83, 8).Value = 0  # H83: 30000 -> 0
$ws.Cells.Item(83, 10).Value = 0  # J83: 30000 -> 0
$ws.Cells.Item(83, 12).Value = 0  # L83: 90000 -> 0
$ws.Cells.Item(83, 14).ClearContents()  # N83: was -99984
$ws.Cells.Item(100, 8).Value = 749  # H100: 695.5333000000001 -> 749
$ws.Cells.Item(100, 9).Value = 627.75  # I100: 657.5454999999999 -> 627.75
$ws.Cells.Item(100, 10).Value = 1040  # J100: 800 -> 1040
$ws.Cells.Item(100, 11).Value = 1255.5  # K100: 1315.091 -> 1255.5
$ws.Cells.Item(100, 12).Value = 2080  # L100: 1600 -> 2080
$ws.Cells.Item(100, 13).Value = -714.5  # M100: -774.0909999999999 -> -714.5
$ws.Cells.Item(100, 14).Value = -3162  # N100: -2682 -> -3162
$ws.Cells.Item(113, 8).Value = 740.9524  # H113: 746.7619 -> 740.9524
$ws.Cells.Item(113, 9).Value = 634.7273  # I113: 600.1667 -> 634.7273
$ws.Cells.Item(113, 10).Value = 857.8  # J113: 942.2222 -> 857.8
$ws.Cells.Item(113, 11).Value = 1904.1819  # K113: 1800.5001 -> 1904.1819
$ws.Cells.Item(113, 12).Value = 2573.4  # L113: 2826.6666 -> 2573.4
$ws.Cells.Item(113, 13).Value = 265.8181  # M113: 369.4999 -> 265.8181
$ws.Cells.Item(113, 14).Value = -6913.4  # N113: -7166.6666 -> -6913.4
$ws.Cells.Item(132, 8).Value = 5729.6665  # H132: 5807.7144 -> 5729.6665
$ws.Cells.Item(132, 9).Value = 5409.033  # I132: 5492.1724 -> 5409.033
$ws.Cells.Item(132, 11).Value = 16227.099  # K132: 16476.5172 -> 16227.099
$ws.Cells.Item(132, 13).Value = -13697.099  # M132: -13946.5172 -> -13697.099
